$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 14497.857
$ws.Range("I12").Value = 14497.857
$ws.Range("J12").Value = 0
$ws.Range("K12").Value = 14497.857
$ws.Range("L12").Value = 0
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -14327.857
$ws.Range("H40").Value = 5311
$ws.Range("I40").Value = 4147.5
$ws.Range("J40").Value = 5776.4
$ws.Range("K40").Value = 4147.5
$ws.Range("L40").Value = 5776.4
$ws.Range("M40").Value = -3972.5
$ws.Range("N40").Value = -6126.4
$ws.Range("H132").Value = 3943.4443
$ws.Range("I132").Value = 3943.4443
$ws.Range("K132").Value = 11830.3329
$ws.Range("M132").Value = -9300.332900000001
$ws.Range("H137").Value = 2088.9473
$ws.Range("J137").Value = 2431.4167
$ws.Range("L137").Value = 7294.250100000001
$ws.Range("N137").Value = -12394.2501

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5662.3447
$ws.Range("I32").Value = 3988.8076
$ws.Range("J32").Value = 20166.334
$ws.Range("K32").Value = 3988.8076
$ws.Range("L32").Value = 20166.334
$ws.Range("M32").Value = -3701.8076
$ws.Range("N32").Value = -20740.334
$ws.Range("H74").Value = 125007700
$ws.Range("I74").Value = 142864370
$ws.Range("K74").Value = 142864370
$ws.Range("M74").Value = -142863496
$ws.Range("H77").Value = 125007700
$ws.Range("I77").Value = 142864370
$ws.Range("K77").Value = 714321850
$ws.Range("M77").Value = -714317482
$ws.Range("H132").Value = 7146688.5
$ws.Range("I132").Value = 9093240
$ws.Range("K132").Value = 27279720
$ws.Range("M132").Value = -27277190

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 425000
$ws.Range("J43").Value = 425000
$ws.Range("L43").Value = 425000
$ws.Range("N43").Value = -425362
$ws.Range("H97").Value = 53719
$ws.Range("I97").Value = 24959
$ws.Range("K97").Value = 24959
$ws.Range("M97").Value = -23968
$ws.Range("H105").Value = 3785.4285
$ws.Range("I105").Value = 2999
$ws.Range("J105").Value = 3916.5
$ws.Range("K105").Value = 2999
$ws.Range("L105").Value = 3916.5
$ws.Range("M105").Value = -1252
$ws.Range("N105").Value = -7410.5
$ws.Range("H107").Value = 101608.8
$ws.Range("I107").Value = 1684.1428
$ws.Range("J107").Value = 334766.34
$ws.Range("K107").Value = 1684.1428
$ws.Range("L107").Value = 334766.34
$ws.Range("M107").Value = 235.8571999999999
$ws.Range("N107").Value = -338606.34
$ws.Range("H115").Value = 112499
$ws.Range("J115").Value = 112499
$ws.Range("L115").Value = 112499
$ws.Range("N115").Value = -115633
$ws.Range("H134").Value = 18524060
$ws.Range("I134").Value = 20838944
$ws.Range("K134").Value = 62516832
$ws.Range("M134").Value = -62514297

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 7049.533
$ws.Range("I22").Value = 10260
$ws.Range("J22").Value = 628.6
$ws.Range("K22").Value = 10260
$ws.Range("L22").Value = 628.6
$ws.Range("M22").Value = -9910
$ws.Range("N22").Value = -1328.6
$ws.Range("H28").Value = 30607.166
$ws.Range("J28").Value = 30607.166
$ws.Range("L28").Value = 30607.166
$ws.Range("N28").Value = -31097.166
$ws.Range("H31").Value = 2980.5186
$ws.Range("I31").Value = 2781.75
$ws.Range("K31").Value = 2781.75
$ws.Range("M31").Value = -2486.75
$ws.Range("H34").Value = 2980.5186
$ws.Range("I34").Value = 2781.75
$ws.Range("K34").Value = 2781.75
$ws.Range("M34").Value = -2579.75
$ws.Range("H58").Value = 17862420
$ws.Range("I58").Value = 38470532
$ws.Range("J58").Value = 2058.2666
$ws.Range("K58").Value = 38470532
$ws.Range("L58").Value = 2058.2666
$ws.Range("M58").Value = -38470329
$ws.Range("N58").Value = -2464.2666
$ws.Range("H62").Value = 2872.25
$ws.Range("I62").Value = 2795
$ws.Range("J62").Value = 2949.5
$ws.Range("K62").Value = 2795
$ws.Range("L62").Value = 2949.5
$ws.Range("M62").Value = -2171
$ws.Range("N62").Value = -4197.5
$ws.Range("H65").Value = 2872.25
$ws.Range("I65").Value = 2795
$ws.Range("J65").Value = 2949.5
$ws.Range("K65").Value = 13975
$ws.Range("L65").Value = 14747.5
$ws.Range("M65").Value = -10855
$ws.Range("N65").Value = -20987.5
$ws.Range("H97").Value = 32666.334
$ws.Range("J97").Value = 32666.334
$ws.Range("L97").Value = 32666.334
$ws.Range("N97").Value = -34648.334
$ws.Range("H122").Value = 2530.818
$ws.Range("I122").Value = 2438.9
$ws.Range("J122").Value = 3450
$ws.Range("K122").Value = 7316.700000000001
$ws.Range("L122").Value = 10350
$ws.Range("M122").Value = -4866.700000000001
$ws.Range("N122").Value = -15250
$ws.Range("H132").Value = 83337360
$ws.Range("I132").Value = 111115070
$ws.Range("J132").Value = 4210.3335
$ws.Range("K132").Value = 333345210
$ws.Range("L132").Value = 12631.0005
$ws.Range("M132").Value = -333342680
$ws.Range("N132").Value = -17691.0005
$ws.Range("H134").Value = 17931306
$ws.Range("I134").Value = 20919358
$ws.Range("J134").Value = 2999.5
$ws.Range("K134").Value = 62758074
$ws.Range("L134").Value = 8998.5
$ws.Range("M134").Value = -62755539
$ws.Range("N134").Value = -14068.5
$ws.Range("H136").Value = 17862420
$ws.Range("I136").Value = 38470532
$ws.Range("J136").Value = 2058.2666
$ws.Range("K136").Value = 115411596
$ws.Range("L136").Value = 6174.7998
$ws.Range("M136").Value = -115409046
$ws.Range("N136").Value = -11274.7998
$ws.Range("H137").Value = 192808.58
$ws.Range("J137").Value = 192808.58
$ws.Range("L137").Value = 192808.58
$ws.Range("N137").Value = -203008.58

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 485.8
$ws.Range("J23").Value = 541
$ws.Range("L23").Value = 1623
$ws.Range("N23").Value = -2093
$ws.Range("H33").Value = 643.2143
$ws.Range("I33").Value = 559.6
$ws.Range("J33").Value = 689.6667
$ws.Range("K33").Value = 3357.6
$ws.Range("L33").Value = 4138.0002
$ws.Range("M33").Value = -3074.6
$ws.Range("N33").Value = -4704.0002
$ws.Range("H118").Value = 0
$ws.Range("I118").Value = 0
$ws.Range("K118").Value = 0
$ws.Range("M118").ClearContents()
$ws.Range("H131").Value = 1631.3846
$ws.Range("J131").Value = 2400
$ws.Range("L131").Value = 7200
$ws.Range("N131").Value = -17280

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 10407.667
$ws.Range("J23").Value = 15556.5
$ws.Range("L23").Value = 15556.5
$ws.Range("N23").Value = -16002.5
$ws.Range("H95").Value = 15000
$ws.Range("J95").Value = 15000
$ws.Range("L95").Value = 15000
$ws.Range("N95").Value = -20492
$ws.Range("H107").Value = 5018.3076
$ws.Range("J107").Value = 7651
$ws.Range("L107").Value = 7651
$ws.Range("N107").Value = -11491
$ws.Range("H132").Value = 6948600
$ws.Range("I132").Value = 7816856.5
$ws.Range("K132").Value = 23450569.5
$ws.Range("M132").Value = -23448039.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 999
$ws.Range("I17").Value = 0
$ws.Range("K17").Value = 0
$ws.Range("M17").ClearContents()
$ws.Range("H132").Value = 21821020
$ws.Range("I132").Value = 22859782
$ws.Range("J132").Value = 7000
$ws.Range("K132").Value = 68579346
$ws.Range("L132").Value = 21000
$ws.Range("M132").Value = -68576816
$ws.Range("N132").Value = -26060
$ws.Range("H136").Value = 1624
$ws.Range("I136").Value = 1488.8182
$ws.Range("K136").Value = 4466.4546
$ws.Range("M136").Value = -1916.4546

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H19").Value = 20005
$ws.Range("I19").Value = 20005
$ws.Range("K19").Value = 20005
$ws.Range("M19").Value = -19831
$ws.Range("H88").Value = 0
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("K88").Value = 0
$ws.Range("L88").ClearContents()
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = 0
$ws.Range("H91").Value = 0
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("K91").Value = 0
$ws.Range("L91").ClearContents()
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = 0
$ws.Range("H94").Value = 18600
$ws.Range("J94").Value = 18600
$ws.Range("L94").Value = 18600
$ws.Range("N94").Value = -20402
$ws.Range("H132").Value = 20009244
$ws.Range("I132").Value = 27784062
$ws.Range("K132").Value = 83352186
$ws.Range("M132").Value = -83349656
$ws.Range("H136").Value = 11364917
$ws.Range("I136").Value = 11364917
$ws.Range("K136").Value = 34094751
$ws.Range("M136").Value = -34092201

Write-Host "Applied all Spriggan_Profits updates"